$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, border, centered) from the existing H1 header
# cell onto the two new header cells so they match the rest of the header row
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the data for the new columns I and J (rows 2-12)
$dataI = @(8, 7, 8, 8, 7, 6, 7, 7, 6, 7, 5)
$dataJ = @(8, 8, 9, 8, 8, 7, 8, 7, 7, 8, 6)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
